# HrHeadAddEmployee test script update (sheet "HRHeadEmployee"):
#   - the username/test-data cell (B3) is changed from the e-mail address
#     "hrhead@gmail.com" to "emil" (the mailto hyperlink that was already on
#     the cell is left untouched)
#   - a password value "pass" is filled into the previously empty D3 cell
# The sheet is then left as the active sheet/selection, matching the saved
# workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HRHeadEmployee")

$ws.Range("D3").Value = "pass"
$ws.Range("B3").Value = "emil"

$ws.Activate()
$ws.Range("B3").Select()
